$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")

# Insert a new row 43 (below the last data row 42), copy row 42's formatting
# down onto it, then clear the Final Score column's inherited format/formula
# since the new row's Final Score is a typed-in plain value, not a formula.
$ws.Rows.Item(43).Insert()
$ws.Range("A42").Copy()
$ws.Range("A43").PasteSpecial(-4122)
$ws.Cells.Item(43, 5).ClearFormats()

# Fill in the new bug-matrix entry.
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "Top k Values"
$ws.Cells.Item(43, 3).Value = "H"
$ws.Cells.Item(43, 4).Value = 5
$ws.Cells.Item(43, 5).Value = 5
$ws.Cells.Item(43, 6).Value = "Top-k Does not account for the last record before midnight properly or there could be some logic issues"
$ws.Cells.Item(43, 7).Value = "N"

$wb.Save()
